$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 61) with the same look & feel as every other
# data row: centered alignment, Date column stored as literal text (not a
# real Excel date), Game column text, ModCount column a plain number.

# 1) Clone the formatting of the last existing data row (row 60) onto the
#    new row 61 in one shot, so we reuse the existing "center/center"
#    style instead of building it up property-by-property (which would
#    create extra transient style entries).
$ws.Range("A60:C60").Copy()
$ws.Range("A61:C61").PasteSpecial(-4122)  # xlPasteFormats

# 2) Column A: "2026/01/10" must land as literal text, exactly like the
#    existing Date cells - not get auto-converted into an Excel date
#    serial number (which is what a plain .Value assignment would do).
#    Stage the exact text in an out-of-the-way scratch cell (forcing text
#    via a leading apostrophe), copy it, then paste only the VALUE into
#    A61 so the center/center formatting already applied above is kept
#    and no stray quote-prefix formatting leaks onto A61 itself.
$scratch = $ws.Range("Z1")
$scratch.Value = "'2026/01/10"
$scratch.Copy()
$ws.Range("A61").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

# 3) Game + ModCount columns.
$ws.Range("B61").Value = "逃离鸭科夫"
$ws.Range("C61").Value = 1141
